# Update "想去人数" (F column) values across the four sheets to the
# newly scraped counts, per the commit "Update gh-pages to output
# generated at 4e2132f".

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 247
$ws1.Range("F5").Value = 1992
$ws1.Range("F6").Value = 79
$ws1.Range("F7").Value = 444
$ws1.Range("F8").Value = 415
$ws1.Range("F9").Value = 214
$ws1.Range("F10").Value = 7002
$ws1.Range("F12").Value = 545
$ws1.Range("F13").Value = 122
$ws1.Range("F15").Value = 2405
$ws1.Range("F16").Value = 1760
$ws1.Range("F17").Value = 146
$ws1.Range("F19").Value = 101
$ws1.Range("F21").Value = 115
$ws1.Range("F23").Value = 172
$ws1.Range("F25").Value = 978
$ws1.Range("F26").Value = 164
$ws1.Range("F27").Value = 4112

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 73
$ws2.Range("F3").Value = 21

# Sheet 3: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 121
$ws3.Range("F3").Value = 696

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 121
$ws4.Range("F4").Value = 696
$ws4.Range("F6").Value = 73
$ws4.Range("F7").Value = 247
$ws4.Range("F8").Value = 1992
$ws4.Range("F9").Value = 21
$ws4.Range("F11").Value = 79
$ws4.Range("F12").Value = 444
$ws4.Range("F13").Value = 415
$ws4.Range("F14").Value = 214
$ws4.Range("F15").Value = 7002
$ws4.Range("F17").Value = 545
$ws4.Range("F18").Value = 122
$ws4.Range("F20").Value = 2405
$ws4.Range("F21").Value = 1760
$ws4.Range("F22").Value = 146
$ws4.Range("F24").Value = 101
$ws4.Range("F26").Value = 115
$ws4.Range("F28").Value = 172
$ws4.Range("F30").Value = 978
$ws4.Range("F31").Value = 164
$ws4.Range("F32").Value = 4112

$wb.Save()
